{"js": "// Applies the Belgian Navy doc edits:\n//   1) \"Anti Submarine\"  -> \"Anti-Submarine\"  (hyphenation fix, also drops the\n//      spell-check \"Anti Submarine\" flag since the hyphenated word is no\n//      longer flagged as a misspelling by the author's edit)\n//   2) \"anti ship\"       -> \"anti-ship\"       (hyphenation fix, same reason)\n//   3) removes the stray \"_GoBack\" bookmark left over from the author's last\n//      editing position (Word drops this automatically on save in recent\n//      builds, which is what the diff shows disappearing)\n//\n// The two text fixes are resolved with Body.search()+insertText() so each\n// edit is scoped to just the matched text (searches are case-sensitive and\n// each phrase is unique in the document, so there is no ambiguity).\n\nconst body = context.document.body;\n\n// 1) \"Anti Submarine\" -> \"Anti-Submarine\"\nconst antiSubmarine = body.search(\"Anti Submarine\", { matchCase: true });\nantiSubmarine.load(\"text\");\nawait context.sync();\n\nif (antiSubmarine.items.length > 0) {\n  antiSubmarine.items[0].insertText(\"Anti-Submarine\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) \"anti ship\" -> \"anti-ship\"\nconst antiShip = body.search(\"anti ship\", { matchCase: true });\nantiShip.load(\"text\");\nawait context.sync();\n\nif (antiShip.items.length > 0) {\n  antiShip.items[0].insertText(\"anti-ship\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3) Remove the leftover \"_GoBack\" bookmark.\nconst goBackRange = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\ngoBackRange.load(\"isNullObject\");\nawait context.sync();\n\nif (!goBackRange.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Applies the Belgian Navy doc edits:\n#   1) \"Anti Submarine\"  -> \"Anti-Submarine\"  (hyphenation fix)\n#   2) \"anti ship\"       -> \"anti-ship\"       (hyphenation fix)\n#   3) removes the stray \"_GoBack\" bookmark left over from the author's last\n#      editing position.\n\n$d = $word.ActiveDocument\n\n# 1) \"Anti Submarine\" -> \"Anti-Submarine\"\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"Anti Submarine\"\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $false\nif ($rng.Find.Execute()) {\n    $rng.Text = \"Anti-Submarine\"\n}\n\n# 2) \"anti ship\" -> \"anti-ship\"\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Text = \"anti ship\"\n$rng2.Find.MatchCase = $true\n$rng2.Find.MatchWholeWord = $false\nif ($rng2.Find.Execute()) {\n    $rng2.Text = \"anti-ship\"\n}\n\n# 3) Remove the leftover \"_GoBack\" bookmark, if present.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n"}
